$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 83
$rowRange = $ws.Range("A$row`:H$row")

# Force text storage (so purely-numeric strings like "246" / phone numbers
# keep their leading context and string type instead of becoming numbers),
# then strip the temporary "Text" number-format back off so the cells end
# up unstyled, matching the rest of the data rows.
$rowRange.NumberFormat = "@"

$ws.Cells.Item($row, 1).Value = "Jabborova Rano Rahmonovna"
$ws.Cells.Item($row, 2).Value = "Maktabgacha talim tashkiloti metodisti"
$ws.Cells.Item($row, 3).Value = "AD2380128"
$ws.Cells.Item($row, 4).Value = "246"
$ws.Cells.Item($row, 5).Value = "Samarqand viloyati"
$ws.Cells.Item($row, 6).Value = "Samarqand tumani"
$ws.Cells.Item($row, 7).Value = "998972861919"
$ws.Cells.Item($row, 8).Value = "13-11-2024"

$rowRange.ClearFormats()
